$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (scale_point = 1): replace the generic "sustainable development goal" text
# with the explicit "lowest potential" need-based prompt.
$ws.Range("B2:I2").Value = "lowest potential for addressing the respective need"

# Row 6 (scale_point = 5): replace the generic "sustainable development goal" text
# with the explicit "highest potential" need-based prompt.
$ws.Range("B6:I6").Value = "highest potential for addressing the respective need"

# Update the active selection to match the new view state.
$ws.Range("G5").Select()
